$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 2969.7778
$ws.Range("I5").Value = 1662.75
$ws.Range("K5").Value = 1662.75
$ws.Range("M5").Value = -1547.75
$ws.Range("H9").Value = 1520.0526
$ws.Range("J9").Value = 1387.8
$ws.Range("L9").Value = 1387.8
$ws.Range("N9").Value = -1725.8
$ws.Range("H15").Value = 103638.24
$ws.Range("I15").Value = 103638.24
$ws.Range("K15").Value = 310914.72
$ws.Range("M15").Value = -310745.72
$ws.Range("H33").Value = 1332.6666
$ws.Range("J33").Value = 1997
$ws.Range("L33").Value = 1997
$ws.Range("N33").Value = -2455
$ws.Range("H41").Value = 708.5
$ws.Range("I41").Value = 318.125
$ws.Range("K41").Value = 318.125
$ws.Range("M41").Value = 121.875
$ws.Range("H58").Value = 5424.5
$ws.Range("J58").Value = 9899.5
$ws.Range("L58").Value = 29698.5
$ws.Range("N58").Value = -29998.5
$ws.Range("H62").Value = 9267341
$ws.Range("I62").Value = 11119582
$ws.Range("K62").Value = 11119582
$ws.Range("M62").Value = -11118958
$ws.Range("H65").Value = 9267341
$ws.Range("I65").Value = 11119582
$ws.Range("K65").Value = 55597910
$ws.Range("M65").Value = -55594790
$ws.Range("H82").Value = 8179.6
$ws.Range("I82").Value = 8999.5
$ws.Range("K82").Value = 26998.5
$ws.Range("M82").Value = -26592.5
$ws.Range("H85").Value = 8179.6
$ws.Range("I85").Value = 8999.5
$ws.Range("K85").Value = 26998.5
$ws.Range("M85").Value = -25594.5
$ws.Range("H96").Value = 250000300
$ws.Range("I96").Value = 398.66666
$ws.Range("K96").Value = 1195.99998
$ws.Range("M96").Value = 177.0000199999999
$ws.Range("H98").Value = 2693.7307
$ws.Range("I98").Value = 2721.48
$ws.Range("K98").Value = 2721.48
$ws.Range("M98").Value = -1223.48
$ws.Range("H99").Value = 100003460
$ws.Range("I99").Value = 885.625
$ws.Range("K99").Value = 2656.875
$ws.Range("M99").Value = -1158.875
$ws.Range("H101").Value = 9995.200000000001
$ws.Range("I101").Value = 9991.375
$ws.Range("K101").Value = 29974.125
$ws.Range("M101").Value = -28352.125
$ws.Range("H122").Value = 2693.7307
$ws.Range("I122").Value = 2721.48
$ws.Range("K122").Value = 8164.440000000001
$ws.Range("M122").Value = -5714.440000000001
$ws.Range("H132").Value = 502573.38
$ws.Range("I132").Value = 579273.2
$ws.Range("K132").Value = 1737819.6
$ws.Range("M132").Value = -1735289.6

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 6125
$ws.Range("I2").Value = 2403.7693
$ws.Range("J2").Value = 15800.2
$ws.Range("K2").Value = 2403.7693
$ws.Range("L2").Value = 15800.2
$ws.Range("M2").Value = -2290.7693
$ws.Range("N2").Value = -16026.2
$ws.Range("H32").Value = 12226937
$ws.Range("I32").Value = 20843602
$ws.Range("J32").Value = 2608333.8
$ws.Range("K32").Value = 20843602
$ws.Range("L32").Value = 2608333.8
$ws.Range("M32").Value = -20843315
$ws.Range("N32").Value = -2608907.8
$ws.Range("H61").Value = 3666.2666
$ws.Range("I61").Value = 3183.5642
$ws.Range("K61").Value = 3183.5642
$ws.Range("M61").Value = -2971.5642
$ws.Range("H97").Value = 2599.4
$ws.Range("I97").Value = 1999.25
$ws.Range("K97").Value = 1999.25
$ws.Range("M97").Value = -1503.25
$ws.Range("H102").Value = 2181.2727
$ws.Range("I102").Value = 2181.2727
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2181.2727
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -559.2727
$ws.Range("N102").ClearContents()
$ws.Range("H116").Value = 6125
$ws.Range("I116").Value = 2403.7693
$ws.Range("J116").Value = 15800.2
$ws.Range("K116").Value = 2403.7693
$ws.Range("L116").Value = 15800.2
$ws.Range("M116").Value = -109.7692999999999
$ws.Range("N116").Value = -20388.2
$ws.Range("H136").Value = 3666.2666
$ws.Range("I136").Value = 3183.5642
$ws.Range("K136").Value = 9550.692599999998
$ws.Range("M136").Value = -7000.692599999998

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 6125
$ws.Range("I3").Value = 2403.7693
$ws.Range("J3").Value = 15800.2
$ws.Range("K3").Value = 2403.7693
$ws.Range("L3").Value = 15800.2
$ws.Range("M3").Value = -2289.7693
$ws.Range("N3").Value = -16028.2
$ws.Range("H86").Value = 7191
$ws.Range("J86").Value = 19172.25
$ws.Range("L86").Value = 19172.25
$ws.Range("N86").Value = -21418.25
$ws.Range("H89").Value = 7191
$ws.Range("J89").Value = 19172.25
$ws.Range("L89").Value = 95861.25
$ws.Range("N89").Value = -107093.25
$ws.Range("H94").Value = 29244.62
$ws.Range("I94").Value = 3978.4849
$ws.Range("K94").Value = 3978.4849
$ws.Range("M94").Value = -3527.4849

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7346
$ws.Range("I31").Value = 1984.375
$ws.Range("K31").Value = 1984.375
$ws.Range("M31").Value = -1689.375
$ws.Range("H34").Value = 7346
$ws.Range("I34").Value = 1984.375
$ws.Range("K34").Value = 1984.375
$ws.Range("M34").Value = -1782.375
$ws.Range("H62").Value = 8281.75
$ws.Range("I62").Value = 8261
$ws.Range("J62").Value = 8323.25
$ws.Range("K62").Value = 8261
$ws.Range("L62").Value = 8323.25
$ws.Range("M62").Value = -7637
$ws.Range("N62").Value = -9571.25
$ws.Range("H65").Value = 8281.75
$ws.Range("I65").Value = 8261
$ws.Range("J65").Value = 8323.25
$ws.Range("K65").Value = 41305
$ws.Range("L65").Value = 41616.25
$ws.Range("M65").Value = -38185
$ws.Range("N65").Value = -47856.25
$ws.Range("H122").Value = 4600.5713
$ws.Range("I122").Value = 4414.067
$ws.Range("J122").Value = 5066.8335
$ws.Range("K122").Value = 13242.201
$ws.Range("L122").Value = 15200.5005
$ws.Range("M122").Value = -10792.201
$ws.Range("N122").Value = -20100.5005

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 1716.7693
$ws.Range("I114").Value = 1371
$ws.Range("J114").Value = 2270
$ws.Range("K114").Value = 4113
$ws.Range("L114").Value = 6810
$ws.Range("M114").Value = -859
$ws.Range("N114").Value = -13318
$ws.Range("H137").Value = 2160.1738
$ws.Range("I137").Value = 2117.4412
$ws.Range("J137").Value = 2281.25
$ws.Range("K137").Value = 6352.323600000001
$ws.Range("L137").Value = 6843.75
$ws.Range("M137").Value = -1252.323600000001
$ws.Range("N137").Value = -17043.75

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 4117.1904
$ws.Range("I93").Value = 3418.9375
$ws.Range("J93").Value = 6351.6
$ws.Range("K93").Value = 3418.9375
$ws.Range("L93").Value = 6351.6
$ws.Range("M93").Value = -2170.9375
$ws.Range("N93").Value = -8847.6

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2584.4285
$ws.Range("I81").Value = 1943.2222
$ws.Range("J81").Value = 3738.6
$ws.Range("K81").Value = 3886.4444
$ws.Range("L81").Value = 7477.2
$ws.Range("M81").Value = -2825.4444
$ws.Range("N81").Value = -9599.200000000001
$ws.Range("H84").Value = 2584.4285
$ws.Range("I84").Value = 1943.2222
$ws.Range("J84").Value = 3738.6
$ws.Range("K84").Value = 19432.222
$ws.Range("L84").Value = 37386
$ws.Range("M84").Value = -14128.222
$ws.Range("N84").Value = -47994
$ws.Range("H136").Value = 8629385
$ws.Range("J136").Value = 11056.8
$ws.Range("L136").Value = 33170.39999999999
$ws.Range("N136").Value = -38270.39999999999
